# Remove the "Figure S3" captioned figure (image) and its image-caption
# paragraph from the end of the appendix, as described by the diff:
# both the CaptionedFigure paragraph (drawing) and the following
# ImageCaption paragraph ("Figure S3. Median posterior predictions ...")
# are deleted, leaving the Figure S2 caption as the last paragraph before
# the section properties.

$d = $word.ActiveDocument

# Locate the ImageCaption paragraph whose text starts with "Figure S3" —
# this is the caption paragraph that follows the figure's drawing
# paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Figure S3")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Figure S3' image caption paragraph"
}

$captionPara = $d.Paragraphs.Item($targetIndex)
$figurePara = $d.Paragraphs.Item($targetIndex - 1)

# Delete the caption paragraph first (so indices/ranges for the figure
# paragraph remain valid), then delete the figure (drawing) paragraph.
$captionPara.Range.Delete()
$figurePara.Range.Delete()
